# Apply targeted edits to slide 9 ("SEC_Manuscript-Figures"):
#   - TextBox 53 ("strong founder effect" label): reposition/resize + fix capitalization
#   - TextBox 54 ("Mean Slope" rotated label): reposition/resize + fix capitalization

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)

# Both target text boxes live inside the single top-level group shape on this slide.
$grp = $s.Shapes.Item(1)

# --- TextBox 53: "strong founder effect" -> "Strong founder effect" ---
$foundBox = $grp.GroupItems.Item(21)

$foundBox.Left = 700.8992125984252
$foundBox.Top = 93.70472440944881
$foundBox.Width = 117.15811023622047
$foundBox.Height = 31.50472440944882

$foundTr = $foundBox.TextFrame.TextRange
$foundPara = $foundTr.Paragraphs(1, 1)
$foundPara.Characters(1, $foundPara.Length).Text = "Strong founder effect"

# --- TextBox 54: "Mean Slope" -> "Mean slope" ---
$slopeBox = $grp.GroupItems.Item(22)

$slopeBox.Left = 163.5637007874016
$slopeBox.Top = 182.05732283464567
$slopeBox.Width = 66.41740157480315
$slopeBox.Height = 19.38748031496063

$slopeTr = $slopeBox.TextFrame.TextRange
$slopeTr.Characters(1, $slopeTr.Length).Text = "Mean slope"
